# Added support for multiple pumps, depending on the pH value.
#
# A new "Dose multiplier interval" column is introduced right after "Dose
# vol." (column G): column H becomes "Dose multiplier interval" and the
# data that used to live in H (and beyond, for row 2 which already carries
# a second pump's columns) shifts one place to the right.
#
# Cell formats (fills/borders) are applied by copying them from an existing
# cell that already carries the right style, via PasteSpecial, so every
# newly-populated cell ends up matching its neighbours instead of staying
# in the default "no style" state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 1: header -------------------------------------------------------
# H1 becomes the new header; I1 keeps the old "Force delay" header text
# (its format already matches, since it is a plain copy of H1's old state).
$ws.Range("I1").Value2 = $ws.Range("H1").Value2
$ws.Range("H1").Value2 = "Dose multiplier interval"

# --- Row 2 (pump 1 block + a second pump block later in the row) --------
# Copy the formats each shifted-into cell needs to end up with, using the
# original, not-yet-touched cells as sources:
#   I2, N2, O2 <- G2 (style used by "Dose vol." / dose-multiplier cells)
#   J2         <- D2 (style used by the following probe-reading cell)
#   L2         <- K2 (style used by the cell after it)
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial($xlPasteFormats)
$ws.Range("N2").PasteSpecial($xlPasteFormats)
$ws.Range("O2").PasteSpecial($xlPasteFormats)

$ws.Range("D2").Copy()
$ws.Range("J2").PasteSpecial($xlPasteFormats)

$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial($xlPasteFormats)

# Now shift the values one column to the right, starting from the end so
# nothing is overwritten before it has been read, then fill in the two new
# "Dose multiplier interval" values (H2 for pump 1, N2 for pump 2).
$ws.Range("O2").Value2 = $ws.Range("M2").Value2
$ws.Range("M2").Value2 = $ws.Range("L2").Value2
$ws.Range("L2").Value2 = $ws.Range("K2").Value2
$ws.Range("K2").Value2 = $ws.Range("J2").Value2
$ws.Range("J2").Value2 = $ws.Range("I2").Value2
$ws.Range("I2").Value2 = $ws.Range("H2").Value2
$ws.Range("H2").Value2 = 0.1
$ws.Range("N2").Value2 = 0.1

# --- Rows 3-6 (single pump blocks) ---------------------------------------
# Each of these rows only had a "Force delay" value in H; that value moves
# to the new I column, and H gets the new dose-multiplier-interval value.
$ws.Range("G3").Copy()
$ws.Range("I3").PasteSpecial($xlPasteFormats)
$ws.Range("I3").Value2 = 1
$ws.Range("H3").Value2 = 0.05

$ws.Range("G4").Copy()
$ws.Range("I4").PasteSpecial($xlPasteFormats)
$ws.Range("I4").Value2 = $ws.Range("H4").Value2
$ws.Range("H4").Value2 = 0.1

$ws.Range("G5").Copy()
$ws.Range("I5").PasteSpecial($xlPasteFormats)
$ws.Range("I5").Value2 = $ws.Range("H5").Value2
$ws.Range("H5").Value2 = 0.1

$ws.Range("G6").Copy()
$ws.Range("I6").PasteSpecial($xlPasteFormats)
$ws.Range("I6").Value2 = $ws.Range("H6").Value2
$ws.Range("H6").Value2 = 0.1

# --- View state: scroll so column D is leftmost, select I10 --------------
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("I10").Select()
